# Fix the "docentes" (teachers) and "estudiantes" (students) sheets so that
# login credentials are generic ("s") and the menu/tableview selection
# state reflects the appropriate active sheet/cell for each user type.

$wb = $excel.ActiveWorkbook

# --- docentes sheet: fix sample login user/password for row 2 ---
$wsDocentes = $wb.Worksheets.Item("docentes")
$wsDocentes.Range("D2").Value = "s"
$wsDocentes.Range("G2").Value = "s"

# --- estudiantes sheet: fix sample login user/password for row 2 ---
$wsEstudiantes = $wb.Worksheets.Item("estudiantes")
$wsEstudiantes.Range("D2").Value = "s"
$wsEstudiantes.Range("G2").Value = "s"

# --- Fix window/selection state: "docentes" becomes the active tab,
#     with H7 selected; "estudiantes" keeps H9 selected but is no longer
#     the active (front) tab. ---
$wsEstudiantes.Activate()
$wsEstudiantes.Range("H9").Select()

$wsDocentes.Activate()
$wsDocentes.Range("H7").Select()
